$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 117.72414  # H33
$ws.Cells.Item(33, 9).Value = 74.083336  # I33
$ws.Cells.Item(33, 11).Value = 74.083336  # K33
$ws.Cells.Item(33, 13).Value = 154.916664  # M33

$ws.Cells.Item(38, 8).Value = 1101.9286  # H38
$ws.Cells.Item(38, 9).Value = 157.55556  # I38
$ws.Cells.Item(38, 10).Value = 2801.8  # J38
$ws.Cells.Item(38, 11).Value = 472.66668  # K38
$ws.Cells.Item(38, 12).Value = 8405.400000000001  # L38
$ws.Cells.Item(38, 13).Value = -100.66668  # M38
$ws.Cells.Item(38, 14).Value = -9149.400000000001  # N38

$ws.Cells.Item(40, 8).Value = 3700.4  # H40
$ws.Cells.Item(40, 9).Value = 4228.857  # I40
$ws.Cells.Item(40, 11).Value = 4228.857  # K40
$ws.Cells.Item(40, 13).Value = -4053.857  # M40

$ws.Cells.Item(86, 8).Value = 2001  # H86
$ws.Cells.Item(86, 9).Value = 1501.5  # I86
$ws.Cells.Item(86, 11).Value = 1501.5  # K86
$ws.Cells.Item(86, 13).Value = -378.5  # M86

$ws.Cells.Item(89, 8).Value = 2001  # H89
$ws.Cells.Item(89, 9).Value = 1501.5  # I89
$ws.Cells.Item(89, 11).Value = 7507.5  # K89
$ws.Cells.Item(89, 13).Value = -1891.5  # M89

$ws.Cells.Item(92, 8).Value = 6536965  # H92
$ws.Cells.Item(92, 9).Value = 8548131  # I92
$ws.Cells.Item(92, 10).Value = 675  # J92
$ws.Cells.Item(92, 11).Value = 8548131  # K92
$ws.Cells.Item(92, 12).Value = 675  # L92
$ws.Cells.Item(92, 13).Value = -8546883  # M92
$ws.Cells.Item(92, 14).Value = -3171  # N92

$ws.Cells.Item(113, 8).Value = 92645.73  # H113
$ws.Cells.Item(113, 9).Value = 101740.4  # I113
$ws.Cells.Item(113, 10).Value = 1699  # J113
$ws.Cells.Item(113, 11).Value = 101740.4  # K113
$ws.Cells.Item(113, 12).Value = 1699  # L113
$ws.Cells.Item(113, 13).Value = -98486.39999999999  # M113
$ws.Cells.Item(113, 14).Value = -8207  # N113

$ws.Cells.Item(132, 8).Value = 199036.06  # H132
$ws.Cells.Item(132, 9).Value = 229814.34  # I132
$ws.Cells.Item(132, 10).Value = 35911.2  # J132
$ws.Cells.Item(132, 11).Value = 689443.02  # K132
$ws.Cells.Item(132, 12).Value = 107733.6  # L132
$ws.Cells.Item(132, 13).Value = -686913.02  # M132
$ws.Cells.Item(132, 14).Value = -112793.6  # N132

$ws.Cells.Item(136, 8).Value = 43500  # H136
$ws.Cells.Item(136, 10).Value = 43500  # J136
$ws.Cells.Item(136, 12).Value = 43500  # L136
$ws.Cells.Item(136, 14).Value = -53700  # N136

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2883.2676  # H32
$ws.Cells.Item(32, 9).Value = 1694.5769  # I32
$ws.Cells.Item(32, 11).Value = 1694.5769  # K32
$ws.Cells.Item(32, 13).Value = -1407.5769  # M32

$ws.Cells.Item(97, 8).Value = 13889370  # H97
$ws.Cells.Item(97, 9).Value = 22222618  # I97
$ws.Cells.Item(97, 10).Value = 624.6667  # J97
$ws.Cells.Item(97, 11).Value = 22222618  # K97
$ws.Cells.Item(97, 12).Value = 624.6667  # L97
$ws.Cells.Item(97, 13).Value = -22222122  # M97
$ws.Cells.Item(97, 14).Value = -1616.6667  # N97

$ws.Cells.Item(139, 8).Value = 75000  # H139
$ws.Cells.Item(139, 10).Value = 75000  # J139
$ws.Cells.Item(139, 12).Value = 75000  # L139
$ws.Cells.Item(139, 14).Value = -85280  # N139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(34, 8).Value = 0  # H34
$ws.Cells.Item(34, 9).Value = 0  # I34
$ws.Cells.Item(34, 11).Value = 0  # K34
$ws.Cells.Item(34, 13).ClearContents()  # M34

$ws.Cells.Item(49, 8).Value = 0  # H49
$ws.Cells.Item(49, 10).Value = 0  # J49
$ws.Cells.Item(49, 12).Value = 0  # L49
$ws.Cells.Item(49, 14).ClearContents()  # N49

$ws.Cells.Item(82, 8).Value = 21944.25  # H82
$ws.Cells.Item(82, 9).Value = 5094.2  # I82
$ws.Cells.Item(82, 10).Value = 50027.668  # J82
$ws.Cells.Item(82, 11).Value = 5094.2  # K82
$ws.Cells.Item(82, 12).Value = 50027.668  # L82
$ws.Cells.Item(82, 13).Value = -4711.2  # M82
$ws.Cells.Item(82, 14).Value = -50793.668  # N82

$ws.Cells.Item(85, 8).Value = 21944.25  # H85
$ws.Cells.Item(85, 9).Value = 5094.2  # I85
$ws.Cells.Item(85, 10).Value = 50027.668  # J85
$ws.Cells.Item(85, 11).Value = 5094.2  # K85
$ws.Cells.Item(85, 12).Value = 50027.668  # L85
$ws.Cells.Item(85, 13).Value = -3768.2  # M85
$ws.Cells.Item(85, 14).Value = -52679.668  # N85

$ws.Cells.Item(112, 8).Value = 0  # H112
$ws.Cells.Item(112, 10).Value = 0  # J112
$ws.Cells.Item(112, 12).Value = 0  # L112
$ws.Cells.Item(112, 14).ClearContents()  # N112

$ws.Cells.Item(118, 8).Value = 0  # H118
$ws.Cells.Item(118, 10).Value = 0  # J118
$ws.Cells.Item(118, 12).Value = 0  # L118
$ws.Cells.Item(118, 14).ClearContents()  # N118

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1202.125  # H31
$ws.Cells.Item(31, 9).Value = 851.2157  # I31
$ws.Cells.Item(31, 10).Value = 1819.2413  # J31
$ws.Cells.Item(31, 11).Value = 851.2157  # K31
$ws.Cells.Item(31, 12).Value = 1819.2413  # L31
$ws.Cells.Item(31, 13).Value = -556.2157  # M31
$ws.Cells.Item(31, 14).Value = -2409.2413  # N31

$ws.Cells.Item(34, 8).Value = 1202.125  # H34
$ws.Cells.Item(34, 9).Value = 851.2157  # I34
$ws.Cells.Item(34, 10).Value = 1819.2413  # J34
$ws.Cells.Item(34, 11).Value = 851.2157  # K34
$ws.Cells.Item(34, 12).Value = 1819.2413  # L34
$ws.Cells.Item(34, 13).Value = -649.2157  # M34
$ws.Cells.Item(34, 14).Value = -2223.2413  # N34

$ws.Cells.Item(99, 8).Value = 6251452.5  # H99
$ws.Cells.Item(99, 9).Value = 12501222  # I99
$ws.Cells.Item(99, 10).Value = 1682.8  # J99
$ws.Cells.Item(99, 11).Value = 12501222  # K99
$ws.Cells.Item(99, 12).Value = 1682.8  # L99
$ws.Cells.Item(99, 13).Value = -12499724  # M99
$ws.Cells.Item(99, 14).Value = -4678.8  # N99

$ws.Cells.Item(122, 8).Value = 1200  # H122
$ws.Cells.Item(122, 9).Value = 0  # I122
$ws.Cells.Item(122, 10).Value = 1200  # J122
$ws.Cells.Item(122, 11).Value = 0  # K122
$ws.Cells.Item(122, 12).Value = 3600  # L122
$ws.Cells.Item(122, 13).ClearContents()  # M122
$ws.Cells.Item(122, 14).Value = -8500  # N122

$ws.Cells.Item(126, 8).Value = 6251452.5  # H126
$ws.Cells.Item(126, 9).Value = 12501222  # I126
$ws.Cells.Item(126, 10).Value = 1682.8  # J126
$ws.Cells.Item(126, 11).Value = 37503666  # K126
$ws.Cells.Item(126, 12).Value = 5048.4  # L126
$ws.Cells.Item(126, 13).Value = -37501196  # M126
$ws.Cells.Item(126, 14).Value = -9988.4  # N126

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1191.4255  # H5
$ws.Cells.Item(5, 9).Value = 734.3103599999999  # I5
$ws.Cells.Item(5, 10).Value = 1927.8889  # J5
$ws.Cells.Item(5, 11).Value = 2202.93108  # K5
$ws.Cells.Item(5, 12).Value = 5783.6667  # L5
$ws.Cells.Item(5, 13).Value = -2090.93108  # M5
$ws.Cells.Item(5, 14).Value = -6007.6667  # N5

$ws.Cells.Item(68, 8).Value = 917.2  # H68
$ws.Cells.Item(68, 9).Value = 720.24243  # I68
$ws.Cells.Item(68, 10).Value = 1299.5294  # J68
$ws.Cells.Item(68, 11).Value = 2160.72729  # K68
$ws.Cells.Item(68, 12).Value = 3898.5882  # L68
$ws.Cells.Item(68, 13).Value = -1349.72729  # M68
$ws.Cells.Item(68, 14).Value = -5520.5882  # N68

$ws.Cells.Item(71, 8).Value = 917.2  # H71
$ws.Cells.Item(71, 9).Value = 720.24243  # I71
$ws.Cells.Item(71, 10).Value = 1299.5294  # J71
$ws.Cells.Item(71, 11).Value = 6482.18187  # K71
$ws.Cells.Item(71, 12).Value = 11695.7646  # L71
$ws.Cells.Item(71, 13).Value = -2426.18187  # M71
$ws.Cells.Item(71, 14).Value = -19807.7646  # N71

$ws.Cells.Item(107, 8).Value = 706.95917  # H107
$ws.Cells.Item(107, 9).Value = 160.79167  # I107
$ws.Cells.Item(107, 10).Value = 1231.28  # J107
$ws.Cells.Item(107, 11).Value = 482.37501  # K107
$ws.Cells.Item(107, 12).Value = 3693.84  # L107
$ws.Cells.Item(107, 13).Value = 1437.62499  # M107
$ws.Cells.Item(107, 14).Value = -7533.84  # N107

$ws.Cells.Item(113, 8).Value = 455.3421  # H113
$ws.Cells.Item(113, 10).Value = 454.91666  # J113
$ws.Cells.Item(113, 12).Value = 1364.74998  # L113
$ws.Cells.Item(113, 14).Value = -5704.749980000001  # N113

$ws.Cells.Item(131, 8).Value = 2395.1047  # H131
$ws.Cells.Item(131, 10).Value = 2614.6882  # J131
$ws.Cells.Item(131, 12).Value = 7844.0646  # L131
$ws.Cells.Item(131, 14).Value = -17924.0646  # N131

$ws.Cells.Item(135, 8).Value = 1191.4255  # H135
$ws.Cells.Item(135, 9).Value = 734.3103599999999  # I135
$ws.Cells.Item(135, 10).Value = 1927.8889  # J135
$ws.Cells.Item(135, 11).Value = 6608.793239999999  # K135
$ws.Cells.Item(135, 12).Value = 17351.0001  # L135
$ws.Cells.Item(135, 13).Value = -4073.793239999999  # M135
$ws.Cells.Item(135, 14).Value = -22421.0001  # N135

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 5000  # H15
$ws.Cells.Item(15, 10).Value = 5000  # J15
$ws.Cells.Item(15, 12).Value = 5000  # L15
$ws.Cells.Item(15, 14).Value = -5576  # N15

$ws.Cells.Item(81, 8).Value = 5000  # H81
$ws.Cells.Item(81, 10).Value = 5000  # J81
$ws.Cells.Item(81, 12).Value = 5000  # L81
$ws.Cells.Item(81, 14).Value = -6996  # N81

$ws.Cells.Item(84, 8).Value = 5000  # H84
$ws.Cells.Item(84, 10).Value = 5000  # J84
$ws.Cells.Item(84, 12).Value = 15000  # L84
$ws.Cells.Item(84, 14).Value = -24984  # N84

$ws.Cells.Item(95, 8).Value = 30344  # H95
$ws.Cells.Item(95, 10).Value = 30344  # J95
$ws.Cells.Item(95, 12).Value = 30344  # L95
$ws.Cells.Item(95, 14).Value = -35836  # N95

$ws.Cells.Item(126, 8).Value = 2301.1538  # H126
$ws.Cells.Item(126, 9).Value = 1901.5  # I126
$ws.Cells.Item(126, 10).Value = 2478.7778  # J126
$ws.Cells.Item(126, 11).Value = 5704.5  # K126
$ws.Cells.Item(126, 12).Value = 7436.3334  # L126
$ws.Cells.Item(126, 13).Value = -3234.5  # M126
$ws.Cells.Item(126, 14).Value = -12376.3334  # N126

$ws.Cells.Item(132, 8).Value = 2516.476  # H132
$ws.Cells.Item(132, 9).Value = 2012.2  # I132
$ws.Cells.Item(132, 10).Value = 5037.857  # J132
$ws.Cells.Item(132, 11).Value = 6036.6  # K132
$ws.Cells.Item(132, 12).Value = 15113.571  # L132
$ws.Cells.Item(132, 13).Value = -3506.6  # M132
$ws.Cells.Item(132, 14).Value = -20173.571  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3855.2  # H40
$ws.Cells.Item(40, 9).Value = 1820.8  # I40
$ws.Cells.Item(40, 10).Value = 4533.3335  # J40
$ws.Cells.Item(40, 11).Value = 1820.8  # K40
$ws.Cells.Item(40, 12).Value = 4533.3335  # L40
$ws.Cells.Item(40, 13).Value = -1684.8  # M40
$ws.Cells.Item(40, 14).Value = -4805.3335  # N40

$ws.Cells.Item(55, 8).Value = 416.83334  # H55
$ws.Cells.Item(55, 10).Value = 651  # J55
$ws.Cells.Item(55, 12).Value = 651  # L55
$ws.Cells.Item(55, 14).Value = -997  # N55

$ws.Cells.Item(80, 8).Value = 0  # H80
$ws.Cells.Item(80, 10).Value = 0  # J80
$ws.Cells.Item(80, 12).Value = 0  # L80
$ws.Cells.Item(80, 14).ClearContents()  # N80

$ws.Cells.Item(83, 8).Value = 0  # H83
$ws.Cells.Item(83, 10).Value = 0  # J83
$ws.Cells.Item(83, 12).Value = 0  # L83
$ws.Cells.Item(83, 14).ClearContents()  # N83

$ws.Cells.Item(96, 8).Value = 0  # H96
$ws.Cells.Item(96, 10).Value = 0  # J96
$ws.Cells.Item(96, 12).Value = 0  # L96
$ws.Cells.Item(96, 14).ClearContents()  # N96

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(99, 8).Value = 34000  # H99
$ws.Cells.Item(99, 10).Value = 34000  # J99
$ws.Cells.Item(99, 12).Value = 34000  # L99
$ws.Cells.Item(99, 14).Value = -39990  # N99

$ws.Cells.Item(132, 8).Value = 15627309  # H132
$ws.Cells.Item(132, 9).Value = 20835044  # I132
$ws.Cells.Item(132, 10).Value = 4106.25  # J132
$ws.Cells.Item(132, 11).Value = 62505132  # K132
$ws.Cells.Item(132, 12).Value = 12318.75  # L132
$ws.Cells.Item(132, 13).Value = -62502602  # M132
$ws.Cells.Item(132, 14).Value = -17378.75  # N132

$ws.Cells.Item(136, 8).Value = 22290694  # H136
$ws.Cells.Item(136, 9).Value = 33433996  # I136
$ws.Cells.Item(136, 11).Value = 100301988  # K136
$ws.Cells.Item(136, 13).Value = -100299438  # M136
